$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Row 1: new base value + recompute chain ---
$ws.Range("B1").Value = 84
$ws.Range("C1").Formula = "=B1+48"
$ws.Range("D1").Formula = "=C1+48"
$ws.Range("E1").Formula = "=D1+48"
$ws.Range("F1:N1").Formula = "=E1+48"

# --- Row 2: new first value + new ratio formulas ---
$ws.Range("A2").Value = 39
$ws.Range("B2").Formula = "=(B1-46)/48"
$ws.Range("C2").Formula = "=(C1-46)/48"
$ws.Range("D2:N2").Formula = "=(D1-46)/48"

# --- Column A rows 3-10 stay driven by the same fill-down formula (values recompute automatically) ---
$ws.Range("A3").Formula = "=A2+48"
$ws.Range("A4:A10").Formula = "=A3+48"

# --- Column B rows 3-10 get the new ratio formula ---
$ws.Range("B3").Formula = "=(A3-21)/48"
$ws.Range("B4:B10").Formula = "=(A4-21)/48"

# --- Style updates: copy the "style 4" look from the already-existing C3/B4 cells ---
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("F7").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Selection moves from C3 to C1 ---
$ws.Range("C1").Select()
